# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.332.32'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.563.98'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = "'1.004"
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = "'289.54"
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = "'0.3756"
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('D8').Value = "'49.13"
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.3383"
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').Value = "'0.07508"
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('D11').Value = "'1.123"
$ws.Range('E11').Value = '  -3.88%  '
$ws.Range('D12').Value = "'1.005"
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = "'20.79"
$ws.Range('E13').Value = '  -3.04%  '
$ws.Range('D14').Value = "'5.896"
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = "'6.858"
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = '1.565.77'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').Value = "'89.42"
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = "'0.06720"
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').Value = "'6.166"
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = "'11.85"
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('D24').Value = '22.361.33'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = "'2.366"
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').Value = "'2.682"
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('D27').Value = "'20.00"
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = "'147.49"
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('D29').Value = "'4.988"
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = "'124.92"
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').Value = '1.739.08'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').Value = "'2.014"
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').Value = "'0.9806"
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = "'5.948"
$ws.Range('E34').Value = '  -4.36%  '
$ws.Range('D35').Value = "'9.824"
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = "'1.407"
$ws.Range('E36').Value = '  +10.15%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').Value = "'0.08457"
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Value = "'0.02450"
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').Value = "'0.2257"
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('D40').Value = "'0.06396"
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').Value = "'5.351"
$ws.Range('E41').Value = '  -3.31%  '
$ws.Range('D42').Value = "'0.6243"
$ws.Range('E42').Value = '  -1.63%  '
$ws.Range('D43').Value = "'10.96"
$ws.Range('E43').Value = '  -6.48%  '
$ws.Range('D44').Value = "'1.004"
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = "'13.80"
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').Value = "'3.789"
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = "'0.5796"
$ws.Range('E47').Value = '  -3.05%  '
$ws.Range('D48').Value = "'2.044"
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('D49').Value = "'1.250"
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('D50').Value = "'123.82"
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = "'0.07303"
$ws.Range('E51').Value = '  +0.34%  '
